$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes the A (Company Name), B (Company Number), H (Category),
# I (SIC Codes), J (SIC Description) and K (Typical Use Case) values across
# rows 2-11 (row 3 is left untouched). Columns C, D, E, F, G are unchanged.
#
# All of these columns are stored as plain text in the workbook (company
# numbers / SIC codes look numeric, and some cells are blank), so every
# write goes through NumberFormat "@" (Text) to stop Excel from silently
# re-interpreting numeric-looking strings as numbers, and the style is
# reset back to "Normal" afterwards so we don't leave a stray text format
# applied to the cell.

function Set-TextCell {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

$cols = @("A", "B", "H", "I", "J", "K")
$rows = @(2, 4, 5, 6, 7, 8, 9, 10, 11)

# Snapshot the current ("before") values for every affected cell first, so
# the later writes (which happen in an arbitrary dictionary order) never
# clobber a value that still needs to be read.
$before = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowData
}

# Mapping: destination row -> source row (i.e. destination row receives the
# values that used to live in the source row).
$mapping = @{
    2  = 4
    4  = 2
    5  = 11
    6  = 10
    7  = 9
    8  = 6
    9  = 5
    10 = 8
    11 = 7
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $before[$srcRow]
    foreach ($c in $cols) {
        Set-TextCell $ws.Range("$c$destRow") $srcData[$c]
    }
}
